$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44320
$ws.Range("N2").Value = 16000
$ws.Range("O2").Value = 17000
$ws.Range("P2").Value = 16500
$ws.Range("S2").Value = 825

# Row 4
$ws.Range("D4").Value = 44798
$ws.Range("M4").Value = 80
$ws.Range("N4").Value = 21000
$ws.Range("O4").Value = 22000
$ws.Range("P4").Value = 21500
$ws.Range("S4").Value = 1075

# Row 5
$ws.Range("D5").Value = 44761
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 20000
$ws.Range("O5").Value = 21000
$ws.Range("P5").Value = 20500
$ws.Range("S5").Value = 1025

# Row 6
$ws.Range("D6").Value = 44708
$ws.Range("M6").Value = 80

# Row 7
$ws.Range("D7").Value = 44533
$ws.Range("M7").Value = 100
$ws.Range("N7").Value = 16000
$ws.Range("O7").Value = 17000
$ws.Range("P7").Value = 16500
$ws.Range("S7").Value = 825

# Row 8
$ws.Range("D8").Value = 44357
$ws.Range("N8").Value = 14000
$ws.Range("O8").Value = 15000
$ws.Range("P8").Value = 14500
$ws.Range("S8").Value = 725
